# Auto-generated edit script: updates crypto price/volume figures
# and swaps the Binance-PegBSC-USD / Bittensor row order (rows 31-32).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.219.39"
$ws.Range("E2").Value = "  +1.16%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.718.46"
$ws.Range("E3").Value = "  +2.20%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.10"
$ws.Range("E5").Value = "  +0.65%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.96"
$ws.Range("E6").Value = "  +1.00%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
# Row 8
$ws.Range("E8").Value = "  -0.62%  "
# Row 9
$ws.Range("E9").Value = "  +6.59%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.10"
$ws.Range("E10").Value = "  +3.55%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.405"
$ws.Range("E11").Value = "  +1.01%  "
# Row 12
$ws.Range("E12").Value = "  +1.52%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "30.42"
$ws.Range("E13").Value = "  +3.98%  "
# Row 14
$ws.Range("E14").Value = "  +15.50%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.208.66"
$ws.Range("E15").Value = "  +2.37%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.036.76"
$ws.Range("E16").Value = "  +1.11%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.727.19"
$ws.Range("E17").Value = "  -4.07%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.87"
$ws.Range("E18").Value = "  +1.57%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.92"
$ws.Range("E19").Value = "  +1.29%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "362.79"
$ws.Range("E20").Value = "  +1.98%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.66"
$ws.Range("E21").Value = "  +4.30%  "
# Row 22
$ws.Range("E22").Value = "  -0.15%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.39"
$ws.Range("E23").Value = "  +2.87%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.79"
$ws.Range("E24").Value = "  +2.26%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000108"
$ws.Range("E25").Value = "  +12.31%  "
# Row 26
$ws.Range("E26").Value = "  -3.57%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.73"
$ws.Range("E27").Value = "  +3.31%  "
# Row 28
$ws.Range("E28").Value = "  +4.45%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.35"
$ws.Range("E29").Value = "  +0.01%  "
# Row 30
$ws.Range("E30").Value = "  +4.38%  "
# Row 31
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "541.25"
$ws.Range("E31").Value = "  -1.40%  "
# Row 32
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.08%  "
# Row 33
$ws.Range("E33").Value = "  -0.14%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.70"
$ws.Range("E34").Value = "  +3.31%  "
# Row 35
$ws.Range("E35").Value = "  -5.28%  "
# Row 36
$ws.Range("E36").Value = "  +1.53%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "21.02"
$ws.Range("E37").Value = "  +3.62%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.65"
$ws.Range("E38").Value = "  -1.20%  "
# Row 39
$ws.Range("E39").Value = "  -2.53%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.05%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "171.86"
$ws.Range("E41").Value = "  +1.92%  "
# Row 42
$ws.Range("E42").Value = "  +0.02%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.97"
$ws.Range("E43").Value = "  +1.41%  "
# Row 44
$ws.Range("E44").Value = "  +1.54%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0620"
$ws.Range("E45").Value = "  +0.99%  "
# Row 46
$ws.Range("E46").Value = "  +2.57%  "
# Row 47
$ws.Range("E47").Value = "  +1.28%  "
# Row 48
$ws.Range("E48").Value = "  +0.86%  "
# Row 49
$ws.Range("E49").Value = "  +5.07%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.69"
$ws.Range("E50").Value = "  +4.63%  "
# Row 51
$ws.Range("E51").Value = "  +0.42%  "
